$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-07-23 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-07-24 Wednesday", 2)

# Update the multiplication problems in the table.
# The table has 20 rows x 5 columns; rows 1, 5, 10, 15, 20 contain the
# actual equations (in document/reading order). Addressing cells
# directly by (row, column) avoids ambiguity from duplicate "239x2="
# text appearing twice in the original document.
$tbl = $d.Tables.Item(1)

$replacements = @(
    @{Row=1;  Col=1; Old="839×7="; New="159×7="},
    @{Row=1;  Col=2; Old="239×2="; New="967×9="},
    @{Row=1;  Col=3; Old="759×6="; New="641×4="},
    @{Row=1;  Col=4; Old="332×9="; New="816×9="},
    @{Row=1;  Col=5; Old="772×5="; New="167×9="},

    @{Row=5;  Col=1; Old="763×8="; New="857×9="},
    @{Row=5;  Col=2; Old="754×5="; New="843×8="},
    @{Row=5;  Col=3; Old="343×8="; New="704×4="},
    @{Row=5;  Col=4; Old="461×7="; New="633×7="},
    @{Row=5;  Col=5; Old="239×2="; New="833×3="},

    @{Row=10; Col=1; Old="511×2="; New="488×6="},
    @{Row=10; Col=2; Old="402×2="; New="818×2="},
    @{Row=10; Col=3; Old="577×5="; New="356×9="},
    @{Row=10; Col=4; Old="393×2="; New="469×6="},
    @{Row=10; Col=5; Old="575×3="; New="432×5="},

    @{Row=15; Col=1; Old="463×8="; New="427×4="},
    @{Row=15; Col=2; Old="380×8="; New="632×9="},
    @{Row=15; Col=3; Old="804×3="; New="224×5="},
    @{Row=15; Col=4; Old="995×2="; New="598×7="},
    @{Row=15; Col=5; Old="557×2="; New="904×5="},

    @{Row=20; Col=1; Old="836×6="; New="288×8="},
    @{Row=20; Col=2; Old="552×3="; New="765×9="},
    @{Row=20; Col=3; Old="592×4="; New="427×2="},
    @{Row=20; Col=4; Old="903×3="; New="780×5="},
    @{Row=20; Col=5; Old="866×4="; New="904×5="}
)

foreach ($r in $replacements) {
    $cell = $tbl.Cell($r.Row, $r.Col)
    $rng = $cell.Range
    # Use wdReplaceOne (1) rather than wdReplaceAll (2) so that the
    # replacement is confined to a single match - this matters because
    # "239x2=" appears twice in the source document.
    $rng.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                       $true, 1, $false, $r.New, 1)
}
